$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3317
$ws.Range("I40").Value = 3242
$ws.Range("K40").Value = 3242
$ws.Range("M40").Value = -3067

$ws.Range("H41").Value = 1979.8572
$ws.Range("I41").Value = 2127.875
$ws.Range("J41").Value = 1782.5
$ws.Range("K41").Value = 2127.875
$ws.Range("L41").Value = 1782.5
$ws.Range("M41").Value = -1687.875
$ws.Range("N41").Value = -2662.5

$ws.Range("H92").Value = 921.5
$ws.Range("I92").Value = 861.44446
$ws.Range("J92").Value = 1101.6666
$ws.Range("K92").Value = 861.44446
$ws.Range("L92").Value = 1101.6666
$ws.Range("M92").Value = 386.55554
$ws.Range("N92").Value = -3597.6666

$ws.Range("H98").Value = 90955220
$ws.Range("I98").Value = 100050700
$ws.Range("K98").Value = 100050700
$ws.Range("M98").Value = -100049202

$ws.Range("H122").Value = 90955220
$ws.Range("I122").Value = 100050700
$ws.Range("K122").Value = 300152100
$ws.Range("M122").Value = -300149650

$ws.Range("H137").Value = 3245.8
$ws.Range("I137").Value = 2309.476
$ws.Range("K137").Value = 6928.428
$ws.Range("M137").Value = -4378.428

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 184.47368
$ws.Range("I5").Value = 301
$ws.Range("K5").Value = 301
$ws.Range("M5").Value = -189

$ws.Range("H32").Value = 9436362
$ws.Range("I32").Value = 10418131
$ws.Range("K32").Value = 10418131
$ws.Range("M32").Value = -10417844

$ws.Range("H45").Value = 2731.6667
$ws.Range("I45").Value = 2411.3333
$ws.Range("K45").Value = 2411.3333
$ws.Range("M45").Value = -2034.3333

$ws.Range("H63").Value = 2990.6667
$ws.Range("I63").Value = 2990.6667
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 2990.6667
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -2304.6667
$ws.Range("N63").ClearContents()

$ws.Range("H66").Value = 2990.6667
$ws.Range("I66").Value = 2990.6667
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 14953.3335
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -11521.3335
$ws.Range("N66").ClearContents()

$ws.Range("H110").Value = 1455.2667
$ws.Range("I110").Value = 1542.3846
$ws.Range("J110").Value = 889
$ws.Range("K110").Value = 1542.3846
$ws.Range("L110").Value = 889
$ws.Range("M110").Value = 502.6153999999999
$ws.Range("N110").Value = -4979

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 184.47368
$ws.Range("I4").Value = 301
$ws.Range("K4").Value = 301
$ws.Range("M4").Value = -186

$ws.Range("H22").Value = 3223.5454
$ws.Range("I22").Value = 5292.8
$ws.Range("K22").Value = 5292.8
$ws.Range("M22").Value = -5119.8

$ws.Range("H80").Value = 1783.6428
$ws.Range("I80").Value = 4036.75
$ws.Range("J80").Value = 882.4
$ws.Range("K80").Value = 4036.75
$ws.Range("L80").Value = 882.4
$ws.Range("M80").Value = -3038.75
$ws.Range("N80").Value = -2878.4

$ws.Range("H83").Value = 1783.6428
$ws.Range("I83").Value = 4036.75
$ws.Range("J83").Value = 882.4
$ws.Range("K83").Value = 20183.75
$ws.Range("L83").Value = 4412
$ws.Range("M83").Value = -15191.75
$ws.Range("N83").Value = -14396

$ws.Range("H105").Value = 932.1667
$ws.Range("I105").Value = 516.4
$ws.Range("J105").Value = 3011
$ws.Range("K105").Value = 516.4
$ws.Range("L105").Value = 3011
$ws.Range("M105").Value = 1230.6
$ws.Range("N105").Value = -6505

$ws.Range("H134").Value = 86287.586
$ws.Range("I134").Value = 2503.3
$ws.Range("K134").Value = 7509.900000000001
$ws.Range("M134").Value = -4974.900000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 349.8
$ws.Range("I7").Value = 549.5
$ws.Range("K7").Value = 549.5
$ws.Range("M7").Value = -436.5

$ws.Range("H31").Value = 716366.6
$ws.Range("I31").Value = 15619.818
$ws.Range("J31").Value = 1066740
$ws.Range("K31").Value = 15619.818
$ws.Range("L31").Value = 1066740
$ws.Range("M31").Value = -15324.818
$ws.Range("N31").Value = -1067330

$ws.Range("H34").Value = 716366.6
$ws.Range("I34").Value = 15619.818
$ws.Range("J34").Value = 1066740
$ws.Range("K34").Value = 15619.818
$ws.Range("L34").Value = 1066740
$ws.Range("M34").Value = -15417.818
$ws.Range("N34").Value = -1067144

$ws.Range("H132").Value = 4700.2666
$ws.Range("I132").Value = 3501.125
$ws.Range("J132").Value = 6070.7144
$ws.Range("K132").Value = 10503.375
$ws.Range("L132").Value = 18212.1432
$ws.Range("M132").Value = -7973.375
$ws.Range("N132").Value = -23272.1432

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2339
$ws.Range("I5").Value = 900.6667
$ws.Range("J5").Value = 4496.5
$ws.Range("K5").Value = 2702.0001
$ws.Range("L5").Value = 13489.5
$ws.Range("M5").Value = -2590.0001
$ws.Range("N5").Value = -13713.5

$ws.Range("H21").Value = 2600
$ws.Range("I21").Value = 4000
$ws.Range("J21").Value = 1200
$ws.Range("K21").Value = 12000
$ws.Range("L21").Value = 3600
$ws.Range("M21").Value = -11827
$ws.Range("N21").Value = -3946

$ws.Range("H76").Value = 4987.5557
$ws.Range("J76").Value = 4987.5557
$ws.Range("L76").Value = 14962.6671
$ws.Range("N76").Value = -15728.6671

$ws.Range("H79").Value = 4987.5557
$ws.Range("J79").Value = 4987.5557
$ws.Range("L79").Value = 14962.6671
$ws.Range("N79").Value = -17614.6671

$ws.Range("H131").Value = 3216.1738
$ws.Range("J131").Value = 3593.9
$ws.Range("L131").Value = 10781.7
$ws.Range("N131").Value = -20861.7

$ws.Range("H135").Value = 2339
$ws.Range("I135").Value = 900.6667
$ws.Range("J135").Value = 4496.5
$ws.Range("K135").Value = 8106.0003
$ws.Range("L135").Value = 40468.5
$ws.Range("M135").Value = -5571.0003
$ws.Range("N135").Value = -45538.5

$ws.Range("H140").Value = 2639.3076
$ws.Range("I140").Value = 2639.3076
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 7917.9228
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = -2737.9228
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 584.1429000000001
$ws.Range("I2").Value = 312
$ws.Range("J2").Value = 788.25
$ws.Range("K2").Value = 312
$ws.Range("L2").Value = 788.25
$ws.Range("M2").Value = -199
$ws.Range("N2").Value = -1014.25

$ws.Range("H102").Value = 3793.3
$ws.Range("I102").Value = 2714.7778
$ws.Range("K102").Value = 2714.7778
$ws.Range("M102").Value = -1092.7778

$ws.Range("H126").Value = 4193.4546
$ws.Range("I126").Value = 4193.222
$ws.Range("J126").Value = 4194.5
$ws.Range("K126").Value = 12579.666
$ws.Range("L126").Value = 12583.5
$ws.Range("M126").Value = -10109.666
$ws.Range("N126").Value = -17523.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 50587.094
$ws.Range("I7").Value = 1796.75
$ws.Range("J7").Value = 115640.89
$ws.Range("K7").Value = 1796.75
$ws.Range("L7").Value = 115640.89
$ws.Range("M7").Value = -1684.75
$ws.Range("N7").Value = -115864.89

$ws.Range("H46").Value = 5184.6313
$ws.Range("I46").Value = 1846.7693
$ws.Range("J46").Value = 12416.667
$ws.Range("K46").Value = 1846.7693
$ws.Range("L46").Value = 12416.667
$ws.Range("M46").Value = -1658.7693
$ws.Range("N46").Value = -12792.667

$ws.Range("H55").Value = 71429176
$ws.Range("J55").Value = 984
$ws.Range("L55").Value = 984
$ws.Range("N55").Value = -1330

$ws.Range("H95").Value = 60000
$ws.Range("J95").Value = 60000
$ws.Range("L95").Value = 60000
$ws.Range("N95").Value = -65492

$ws.Range("H126").Value = 50587.094
$ws.Range("I126").Value = 1796.75
$ws.Range("J126").Value = 115640.89
$ws.Range("K126").Value = 5390.25
$ws.Range("L126").Value = 346922.67
$ws.Range("M126").Value = -2920.25
$ws.Range("N126").Value = -351862.67

$ws.Range("H136").Value = 69404.78999999999
$ws.Range("J136").Value = 109909.18
$ws.Range("L136").Value = 329727.54
$ws.Range("N136").Value = -334827.54

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8703606
$ws.Range("I62").Value = 8345.5
$ws.Range("J62").Value = 40006544
$ws.Range("K62").Value = 8345.5
$ws.Range("L62").Value = 40006544
$ws.Range("M62").Value = -7721.5
$ws.Range("N62").Value = -40007792

$ws.Range("H65").Value = 8703606
$ws.Range("I65").Value = 8345.5
$ws.Range("J65").Value = 40006544
$ws.Range("K65").Value = 41727.5
$ws.Range("L65").Value = 200032720
$ws.Range("M65").Value = -38607.5
$ws.Range("N65").Value = -200038960

$ws.Range("H81").Value = 6481.5
$ws.Range("I81").Value = 953.75
$ws.Range("K81").Value = 1907.5
$ws.Range("M81").Value = -846.5

$ws.Range("H84").Value = 6481.5
$ws.Range("I84").Value = 953.75
$ws.Range("K84").Value = 9537.5
$ws.Range("M84").Value = -4233.5

$ws.Range("H97").Value = 26524
$ws.Range("J97").Value = 26524
$ws.Range("L97").Value = 26524
$ws.Range("N97").Value = -28506

$ws.Range("H100").Value = 1456.0834
$ws.Range("I100").Value = 1461.1818
$ws.Range("K100").Value = 2922.3636
$ws.Range("M100").Value = -2381.3636

$ws.Range("H132").Value = 3666.6667
$ws.Range("I132").Value = 3666.6667
$ws.Range("K132").Value = 11000.0001
$ws.Range("M132").Value = -8470.0001
